{"js": "// Replace every occurrence of the old \"\u041f\u0435\u0433\u0430\u0437 2022:\" sentence with the\n// updated wording that adds \"\u0442\u043e\u043a\u043e\u043c 2022. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e\" before the\n// date ranges. The same run appears four times in this document (once\n// per language/section variant), so search the whole body and update\n// every match found.\nconst oldText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0441\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0433\u0430\u0437 2022: 8-17. \u043e\u043a\u0442\u043e\u0431\u0430\u0440, 7-16. \u043d\u043e\u0432\u0435\u043c\u0431\u0430\u0440,\";\nconst newText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0441\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0433\u0430\u0437 \u0442\u043e\u043a\u043e\u043c 2022. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 8-17. \u043e\u043a\u0442\u043e\u0431\u0430\u0440, 7-16. \u043d\u043e\u0432\u0435\u043c\u0431\u0430\u0440,\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the constellation announcement sentence wherever it occurs in the\n# document: four paragraphs (one per language/section variant) contain the\n# exact same run \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0441\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0433\u0430\u0437 2022: 8-17. \u043e\u043a\u0442\u043e\u0431\u0430\u0440, 7-16.\n# \u043d\u043e\u0432\u0435\u043c\u0431\u0430\u0440,\". Find & Replace with \"Replace All\" updates every instance in\n# one pass.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0441\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0433\u0430\u0437 2022: 8-17. \u043e\u043a\u0442\u043e\u0431\u0430\u0440, 7-16. \u043d\u043e\u0432\u0435\u043c\u0431\u0430\u0440,\"\n$find.Replacement.Text = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0441\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0433\u0430\u0437 \u0442\u043e\u043a\u043e\u043c 2022. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 8-17. \u043e\u043a\u0442\u043e\u0431\u0430\u0440, 7-16. \u043d\u043e\u0432\u0435\u043c\u0431\u0430\u0440,\"\n\n$find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n"}
